$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 8, leaving only header (row 1) + data row (row 2)
$ws.Rows("3:8").Delete()

# Fix the path in B2: the first backslash after "обрізані фото" becomes a forward slash
$ws.Range("B2").Value = "C:/Users/Asus/Desktop/не всі фото/обрізані фото/BN-GC-14-1-o-felt-d\2_image_BN-GC-14-1-o-felt-d-Photoroom.jpg"
